# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row (B=TB, C=d2S, D=K, E=IP, F=Win, G=sum)
$data = @{
    2 = @(3.182878228561681, 1.65323645889881, 16.98373111632243, 0.4998867070740569, 0, 22.31973251085698)
    3 = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 0, 6.048734245549538)
    4 = @(0.1554434735375247, 0.3375848360084654, 0.1529057820181812, 0.4998867070740569, 0, 1.145820798638228)
    5 = @(3.182878228561681, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 0, 5.488907176552729)
    6 = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 1, 6.048734245549538)
    7 = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 0, 6.048734245549538)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]  # B
    $ws.Cells.Item($row, 3).Value = $vals[1]  # C
    $ws.Cells.Item($row, 4).Value = $vals[2]  # D
    $ws.Cells.Item($row, 5).Value = $vals[3]  # E
    $ws.Cells.Item($row, 6).Value = $vals[4]  # F
    $ws.Cells.Item($row, 7).Value = $vals[5]  # G
}

$wb.Save()
